$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 1-19 with new filtered GPS values
$ws.Range("A1").Value = -8.309860229492188
$ws.Range("B1").Value = 25.91147804260254
$ws.Range("A2").Value = -8.309860229492188
$ws.Range("B2").Value = 25.91147804260254
$ws.Range("A3").Value = -7.923340797424316
$ws.Range("B3").Value = 27.88959503173828
$ws.Range("A4").Value = -7.506770610809326
$ws.Range("B4").Value = 30.1965274810791
$ws.Range("A5").Value = -7.200259208679199
$ws.Range("B5").Value = 32.2003288269043
$ws.Range("A6").Value = -6.954334735870361
$ws.Range("B6").Value = 34.19545745849609
$ws.Range("A7").Value = -6.746068477630615
$ws.Range("B7").Value = 36.48851013183594
$ws.Range("A8").Value = -6.617000579833984
$ws.Range("B8").Value = 38.66115570068359
$ws.Range("A9").Value = -6.540402889251709
$ws.Range("B9").Value = 40.83989334106445
$ws.Range("A10").Value = -6.50297212600708
$ws.Range("B10").Value = 42.98240280151367
$ws.Range("A11").Value = -6.492049217224121
$ws.Range("B11").Value = 45.23105239868164
$ws.Range("A12").Value = -6.50359058380127
$ws.Range("B12").Value = 47.43316268920898
$ws.Range("A13").Value = -6.53030252456665
$ws.Range("B13").Value = 49.46305084228516
$ws.Range("A14").Value = -6.565435409545898
$ws.Range("B14").Value = 51.55378341674805
$ws.Range("A15").Value = -6.597684860229492
$ws.Range("B15").Value = 53.56374740600586
$ws.Range("A16").Value = -6.620381832122803
$ws.Range("B16").Value = 55.65523147583008
$ws.Range("A17").Value = -6.629554748535156
$ws.Range("B17").Value = 57.67204284667969
$ws.Range("A18").Value = -6.626868724822998
$ws.Range("B18").Value = 59.97574234008789
$ws.Range("A19").Value = -6.616115093231201
$ws.Range("B19").Value = 61.9943962097168

# Remove the remaining rows (20-81) that are no longer needed
$ws.Rows("20:81").Delete()

Write-Output "Done"
